$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- capture style templates before clearing (placed far outside final A1:Y30 dimension) ---
$ws.Range("A1").Copy()
$ws.Range("AA1").PasteSpecial(-4122)   # style s=4 (title)
$ws.Range("I4").Copy()
$ws.Range("AA2").PasteSpecial(-4122)   # style s=2 (yellow marker)
$ws.Range("H20").Copy()
$ws.Range("AA3").PasteSpecial(-4122)   # style s=3 (theme5 marker)
$ws.Range("C3").Copy()
$ws.Range("AA4").PasteSpecial(-4122)   # style s=5 (green marker)
$excel.CutCopyMode = $false

# --- clear old data (row 2 header numbers are unchanged, left untouched) ---
$ws.Range("A1:Y1").Clear()
$ws.Range("A3:Y36").Clear()

# --- row 1 title ---
$ws.Range("A1").Value = "PP4 Project plan"
$ws.Range("AA1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 4 ---
$ws.Range("A4").Value = "Blog"
$ws.Range("B4").Value = "Refine comments and detail views"
$ws.Range("AA4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 5 ---
$ws.Range("A5").Value = "Blog"
$ws.Range("B5").Value = "Add likes functionality"
$ws.Range("AA4").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 6 ---
$ws.Range("A6").Value = "Blog"
$ws.Range("B6").Value = "Link recent blogs to homepage"
$ws.Range("AA4").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 7 ---
$ws.Range("A7").Value = "Blog"
$ws.Range("B7").Value = "Tidy Blog views"
$ws.Range("AA4").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 8 ---
$ws.Range("A8").Value = "General"
$ws.Range("B8").Value = "add contact me link to hompage"
$ws.Range("AA4").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 9 ---
$ws.Range("A9").Value = "General"
$ws.Range("B9").Value = "add info section to homepage"
$ws.Range("AA4").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 10 ---
$ws.Range("A10").Value = "General"
$ws.Range("B10").Value = "add social links to footer"
$ws.Range("AA4").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 11 ---
$ws.Range("A11").Value = "Contact me"
$ws.Range("B11").Value = "Design simple contact form with map and details"
$ws.Range("AA4").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 12 ---
$ws.Range("A12").Value = "Contact me"
$ws.Range("B12").Value = "Add verification"
$ws.Range("AA4").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 13 ---
$ws.Range("A13").Value = "Contact me"
$ws.Range("B13").Value = "email form and show confirmation message"
$ws.Range("AA4").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 14 ---
$ws.Range("A14").Value = "Book"
$ws.Range("B14").Value = "Create model"
$ws.Range("AA2").Copy()
$ws.Range("J14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 15 ---
$ws.Range("A15").Value = "Book"
$ws.Range("B15").Value = "Create view that shows 4 options"
$ws.Range("AA2").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 16 ---
$ws.Range("A16").Value = "Book"
$ws.Range("B16").Value = "Add button to book sessions"
$ws.Range("AA2").Copy()
$ws.Range("K16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 17 ---
$ws.Range("A17").Value = "Book"
$ws.Range("B17").Value = "Select date and time for session and confirm"
$ws.Range("AA2").Copy()
$ws.Range("K17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 18 ---
$ws.Range("A18").Value = "General"
$ws.Range("B18").Value = "User panel creation"
$ws.Range("AA2").Copy()
$ws.Range("N18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 19 ---
$ws.Range("A19").Value = "General"
$ws.Range("B19").Value = "Sam to design Images"
$ws.Range("AA3").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("AA3").Copy()
$ws.Range("I19").PasteSpecial(-4122)
$ws.Range("AA3").Copy()
$ws.Range("J19").PasteSpecial(-4122)
$ws.Range("AA3").Copy()
$ws.Range("K19").PasteSpecial(-4122)
$ws.Range("AA3").Copy()
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("AA3").Copy()
$ws.Range("O19").PasteSpecial(-4122)
$ws.Range("AA3").Copy()
$ws.Range("P19").PasteSpecial(-4122)
$ws.Range("AA3").Copy()
$ws.Range("Q19").PasteSpecial(-4122)
$ws.Range("AA3").Copy()
$ws.Range("R19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 20 ---
$ws.Range("A20").Value = "Book"
$ws.Range("B20").Value = "Add booking to user panel"
$ws.Range("AA2").Copy()
$ws.Range("O20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 21 ---
$ws.Range("A21").Value = "Mobile"
$ws.Range("B21").Value = "Fix any issues with responsivness"
$ws.Range("AA2").Copy()
$ws.Range("P21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 22 ---
$ws.Range("A22").Value = "Test"
$ws.Range("B22").Value = "Create tests for views and models as required"
$ws.Range("AA2").Copy()
$ws.Range("Q22").PasteSpecial(-4122)
$ws.Range("AA2").Copy()
$ws.Range("R22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 23 ---
$ws.Range("A23").Value = "Test"
$ws.Range("B23").Value = "use usual online tools to verify code"
$ws.Range("AA2").Copy()
$ws.Range("Q23").PasteSpecial(-4122)
$ws.Range("AA2").Copy()
$ws.Range("R23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 24 ---
$ws.Range("A24").Value = "Readme"
$ws.Range("B24").Value = "Create readme"
$ws.Range("AA2").Copy()
$ws.Range("U24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 25 ---
$ws.Range("A25").Value = "User Stories"
$ws.Range("B25").Value = "Blog"
$ws.Range("AA2").Copy()
$ws.Range("U25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 26 ---
$ws.Range("A26").Value = "User Stories"
$ws.Range("B26").Value = "Booking"
$ws.Range("AA2").Copy()
$ws.Range("U26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 27 ---
$ws.Range("A27").Value = "User Stories"
$ws.Range("B27").Value = "Contact"
$ws.Range("AA2").Copy()
$ws.Range("U27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 28 ---
$ws.Range("A28").Value = "User Stories"
$ws.Range("B28").Value = "General"
$ws.Range("AA2").Copy()
$ws.Range("U28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 29 ---
$ws.Range("A29").Value = "User Stories"
$ws.Range("B29").Value = "User Panel"
$ws.Range("AA2").Copy()
$ws.Range("U29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- row 30 ---
$ws.Range("A30").Value = "Docs"
$ws.Range("B30").Value = "Check project docs to confirm all areas are covered"
$ws.Range("AA2").Copy()
$ws.Range("U30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- remove scratch template cells so they do not affect the used range ---
$ws.Range("AA1:AA4").Clear()

# --- final selection state (matches author saving with N19 selected) ---
$ws.Range("N19").Select()

Write-Output "done"